$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> updated column values (only columns that change)
$updates = @{
    4  = @{ E = 16 }
    15 = @{ E = 133; F = 66; H = 66 }
    16 = @{ E = 8 }
    17 = @{ E = 86 }
    18 = @{ E = 84; F = 32; H = 32 }
    19 = @{ E = 36 }
    20 = @{ E = 4 }
    24 = @{ E = 18 }
    29 = @{ E = 13 }
    31 = @{ E = 2; F = 1; H = 1 }
    32 = @{ E = 14 }
    33 = @{ E = 26 }
    36 = @{ E = 70; F = 26; H = 26 }
    37 = @{ E = 35 }
    38 = @{ E = 53 }
    39 = @{ E = 18 }
    41 = @{ E = 25 }
    42 = @{ E = 28; F = 11; H = 11 }
    46 = @{ E = 20 }
    47 = @{ E = 45 }
    49 = @{ E = 52 }
    50 = @{ E = 17 }
    61 = @{ E = 21 }
    62 = @{ E = 32; F = 6; H = 6 }
    63 = @{ E = 19 }
    65 = @{ E = 26 }
    72 = @{ E = 28 }
    75 = @{ E = 11 }
    77 = @{ E = 40 }
    78 = @{ E = 33 }
    79 = @{ E = 23 }
    81 = @{ E = 10 }
    82 = @{ E = 10 }
    88 = @{ E = 17 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
